$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.11") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '95.681.99'
$ws.Range("E2").Value = '  +4.02%  '
$ws.Range("D3").Value = '3.088.53'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '235.86'
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("D6").Value = '603.28'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").Value = '1.11'
$ws.Range("E7").Value = '  +2.81%  '
$ws.Range("D8").Value = '0.379'
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D10").Value = '3.081.58'
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").Value = '0.782'
$ws.Range("E11").Value = '  +2.17%  '
$ws.Range("D12").Value = '0.196'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '94.956.83'
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").Value = '0.0000236'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '33.22'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").Value = '5.31'
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("D17").Value = '3.655.30'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '3.062.30'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").Value = '3.51'
$ws.Range("E19").Value = '  -6.79%  '
$ws.Range("D20").Value = '14.23'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").Value = '447.98'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").Value = '5.58'
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").Value = '0.0000189'
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").Value = '8.65'
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("D25").Value = '5.47'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").Value = '84.50'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '11.52'
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").Value = '3.233.80'
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").Value = '0.130'
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("D31").Value = '0.241'
$ws.Range("E31").Value = '  +4.64%  '
$ws.Range("D32").Value = '0.177'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -4.06%  '
$ws.Range("D34").Value = '8.87'
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("D35").Value = '25.44'
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '7.27'
$ws.Range("E36").Value = '  -6.53%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.149'
$ws.Range("E37").Value = '  -4.31%  '
$ws.Range("D38").Value = '483.47'
$ws.Range("E38").Value = '  +3.74%  '
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '24.10'
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '3.69'
$ws.Range("E41").Value = '  -4.56%  '
$ws.Range("D42").Value = '0.429'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("E43").Value = '  -3.83%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '3.12'
$ws.Range("E45").Value = '  -3.91%  '
$ws.Range("D46").Value = '160.71'
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = '0.672'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '1.80'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = '0.000276'
$ws.Range("E49").Value = '  +14.28%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '43.71'
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.12%  '
